# The deck originally ships two DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" (stock Office palette)
#   ppt/theme/theme2.xml -> "Integral"     (the palette actually driving
#                                            the slide master / slides)
# The authored change swaps the two themes' content so the slide master
# (and therefore every slide) now renders with the stock "Office Theme"
# colour scheme instead of "Integral".
#
# Convert an RRGGBB hex string into the little-endian integer that the
# PowerPoint object model's RGB properties expect (0x00BBGGRR).
function HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# The 12 theme colour scheme slots, in their canonical OOXML order.
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToRgb $officeThemeColors[$i - 1]
}
